# Test020.2 / Test021 update
# - Move workbook window position
# - Update TST001 sheet: fix typo, change "Public" study flow to "My Studies / QA3 Group" flow
# - Re-date E2
# - Recolor a few edited cells (new blue font) and move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TST001")

# --- Window position (bookViews/workbookView xWindow/yWindow) ---
$win = $excel.ActiveWindow
$win.Left = -105
$win.Top = 75

# --- Row 2: re-date the run ---
$d = Get-Date -Year 2011 -Month 5 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = $d

# --- Row 4: fix "Wecome" typo ---
$ws.Range("H4").Value = 'Welcome to caIntegrator page loaded'

# --- Rows 9-10: Public Studies -> My Studies wording ---
$ws.Range("G9").Value = 'Verify that in My Studies drop-down menu, the following studies are present: qa1adminpub, qa2adminpub, qa1managpub, qa2managpub'
$ws.Range("G10").Value = 'Verify that in My Studies drop-down menu, the following studies are not present: qa1adminpri, qa2adminpri, qa1managpri, qa2managpri'

# --- Row 11: select QA3 Group study instead of Public study ---
$ws.Range("G11").Value = 'In My Studies drop-down menu a the top, select "QA Test Study - QA3 Group"'
$ws.Range("H11").Value = 'Welcome to QA Test Study -QA3 Group page loaded'
$ws.Range("G11").Font.Color = 12611584
$ws.Range("H11").Font.Color = 12611584

# --- Row 17: welcome text updated to match QA3 Group study ---
$ws.Range("H17").Value = 'Welcome to QA Test Study -QA3 Group page loaded'
$ws.Range("H17").Font.Color = 12611584

# --- Row 18: search QA3 Group study instead of Public study ---
$ws.Range("G18").Value = 'Click on Search QA Test Study - QA3 Group in left menu'
$ws.Range("H18").Value = 'Search QA Test Study - QA3 Group page loaded'
$ws.Range("G18").Font.Color = 12611584
$ws.Range("H18").Font.Color = 12611584

# --- Move the active selection ---
$ws.Range("G14").Select()
